$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45177 = 2023-09-08)
# that was bumped by one day (45178 = 2023-09-09) for every data row (2..439).
for ($r = 2; $r -le 439; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45177) {
        $cell.Value = 45178
    }
}
